$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Actualiza base de datos EC: corrige Salario Basico y Valor Mora
$ws.Range("G16").Value = 877803
$ws.Range("F17").Value = 60000
$ws.Range("F23").Value = 22000
